$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (row 7) first so the new shared strings are appended
# to the shared-string table before the corrected header/name strings.
$ws.Range("A7").Value = "New1"
$ws.Range("B7").Value = "New2"
$ws.Range("C7").Value = "New3"

# Fix casing/typos of existing header and name cells.
$ws.Range("C1").Value = "Header3"
$ws.Range("A1").Value = "Header1"
$ws.Range("C3").Value = "Name3"

# Move the active selection to C3.
[void]$ws.Range("C3").Select()
